# Case_4_27 vm_pu.xlsx: update bus voltage-magnitude results for the 380 kV case
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.035608601327538
$ws.Cells.Item(2, 4).Value = 1.0357846216343
$ws.Cells.Item(2, 5).Value = 1.048839284058674
$ws.Cells.Item(2, 6).Value = 1.05567303597649
$ws.Cells.Item(2, 9).Value = 1.031675728340551
$ws.Cells.Item(2, 10).Value = 1.040721444858092
$ws.Cells.Item(2, 11).Value = 1.038580233683906
$ws.Cells.Item(2, 12).Value = 1.051598010514081
$ws.Cells.Item(2, 13).Value = 1.05841284808516
$ws.Cells.Item(2, 14).Value = 1.042199388682209
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.036729987729014
$ws.Cells.Item(3, 4).Value = 1.03657314809691
$ws.Cells.Item(3, 5).Value = 1.049984095778473
$ws.Cells.Item(3, 6).Value = 1.05696197263096
$ws.Cells.Item(3, 9).Value = 1.031835198188935
$ws.Cells.Item(3, 10).Value = 1.041485583483481
$ws.Cells.Item(3, 11).Value = 1.039178282297128
$ws.Cells.Item(3, 12).Value = 1.052554088126173
$ws.Cells.Item(3, 13).Value = 1.059514060685672
$ws.Cells.Item(3, 14).Value = 1.042964612472093
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.037455588647593
$ws.Cells.Item(4, 4).Value = 1.037083086987379
$ws.Cells.Item(4, 5).Value = 1.050725266347014
$ws.Cells.Item(4, 6).Value = 1.057796679758567
$ws.Cells.Item(4, 9).Value = 1.03193689003659
$ws.Cells.Item(4, 10).Value = 1.041979484220437
$ws.Cells.Item(4, 11).Value = 1.039564315751494
$ws.Cells.Item(4, 12).Value = 1.053172553978173
$ws.Cells.Item(4, 13).Value = 1.060226744445264
$ws.Cells.Item(4, 14).Value = 1.043459214604747
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.037760629617176
$ws.Cells.Item(5, 4).Value = 1.037297395544823
$ws.Cells.Item(5, 5).Value = 1.051036951117942
$ws.Cells.Item(5, 6).Value = 1.058147754246698
$ws.Cells.Item(5, 9).Value = 1.031979283064401
$ws.Cells.Item(5, 10).Value = 1.042186989259539
$ws.Cells.Item(5, 11).Value = 1.039726378239741
$ws.Cells.Item(5, 12).Value = 1.053432514457256
$ws.Cells.Item(5, 13).Value = 1.060526387596745
$ws.Cells.Item(5, 14).Value = 1.043667014324806
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.037811847269996
$ws.Cells.Item(6, 4).Value = 1.037333374797994
$ws.Cells.Item(6, 5).Value = 1.051089290051361
$ws.Cells.Item(6, 6).Value = 1.058206710823671
$ws.Cells.Item(6, 9).Value = 1.031986380030167
$ws.Cells.Item(6, 10).Value = 1.042221822607556
$ws.Cells.Item(6, 11).Value = 1.039753575993748
$ws.Cells.Item(6, 12).Value = 1.053476160486935
$ws.Cells.Item(6, 13).Value = 1.060576700836677
$ws.Cells.Item(6, 14).Value = 1.043701897140173
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.037459664625222
$ws.Cells.Item(7, 4).Value = 1.037085950863175
$ws.Cells.Item(7, 5).Value = 1.050729430713476
$ws.Cells.Item(7, 6).Value = 1.05780137018883
$ws.Cells.Item(7, 9).Value = 1.031937457902041
$ws.Cells.Item(7, 10).Value = 1.041982257426539
$ws.Cells.Item(7, 11).Value = 1.039566482126315
$ws.Cells.Item(7, 12).Value = 1.053176027749326
$ws.Cells.Item(7, 13).Value = 1.060230748169263
$ws.Cells.Item(7, 14).Value = 1.043461991749121
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.035987582001378
$ws.Cells.Item(8, 4).Value = 1.036051167857035
$ws.Cells.Item(8, 5).Value = 1.049226095967267
$ws.Cells.Item(8, 6).Value = 1.056108498225212
$ws.Cells.Item(8, 9).Value = 1.031729931564167
$ws.Cells.Item(8, 10).Value = 1.040979802786899
$ws.Cells.Item(8, 11).Value = 1.038782542482955
$ws.Cells.Item(8, 12).Value = 1.051921159110348
$ws.Cells.Item(8, 13).Value = 1.058784982542586
$ws.Cells.Item(8, 14).Value = 1.042458113508909
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.033393454805424
$ws.Cells.Item(9, 4).Value = 1.03422553234103
$ws.Cells.Item(9, 5).Value = 1.046580066085211
$ws.Cells.Item(9, 6).Value = 1.053130592813287
$ws.Cells.Item(9, 9).Value = 1.031352792860832
$ws.Cells.Item(9, 10).Value = 1.039209135068871
$ws.Cells.Item(9, 11).Value = 1.037393909546728
$ws.Cells.Item(9, 12).Value = 1.049708512246944
$ws.Cells.Item(9, 13).Value = 1.056238282704161
$ws.Cells.Item(9, 14).Value = 1.040684931239622
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.03166388446113
$ws.Cells.Item(10, 4).Value = 1.033006961645809
$ws.Cells.Item(10, 5).Value = 1.044818041365609
$ws.Cells.Item(10, 6).Value = 1.051148715529029
$ws.Cells.Item(10, 9).Value = 1.031093674323592
$ws.Cells.Item(10, 10).Value = 1.038025829933611
$ws.Cells.Item(10, 11).Value = 1.036463288188596
$ws.Cells.Item(10, 12).Value = 1.048232427164639
$ws.Cells.Item(10, 13).Value = 1.054541048543372
$ws.Cells.Item(10, 14).Value = 1.039499945675344
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.030914910007909
$ws.Cells.Item(11, 4).Value = 1.032478957055321
$ws.Cells.Item(11, 5).Value = 1.044055527021875
$ws.Cells.Item(11, 6).Value = 1.050291329038869
$ws.Cells.Item(11, 9).Value = 1.030979649456186
$ws.Cells.Item(11, 10).Value = 1.03751276077853
$ws.Cells.Item(11, 11).Value = 1.036059162972217
$ws.Cells.Item(11, 12).Value = 1.047593022059845
$ws.Cells.Item(11, 13).Value = 1.053806248437817
$ws.Cells.Item(11, 14).Value = 1.038986147903212
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.030636697252763
$ws.Cells.Item(12, 4).Value = 1.032282779224189
$ws.Cells.Item(12, 5).Value = 1.043772362363964
$ws.Cells.Item(12, 6).Value = 1.049972973102784
$ws.Cells.Item(12, 9).Value = 1.030937021400315
$ws.Cells.Item(12, 10).Value = 1.037322079687689
$ws.Cells.Item(12, 11).Value = 1.035908878346084
$ws.Cells.Item(12, 12).Value = 1.047355479950667
$ws.Cells.Item(12, 13).Value = 1.053533326777911
$ws.Cells.Item(12, 14).Value = 1.038795196023351
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.030696375321824
$ws.Cells.Item(13, 4).Value = 1.032324862452439
$ws.Cells.Item(13, 5).Value = 1.043833099124432
$ws.Cells.Item(13, 6).Value = 1.050041256325905
$ws.Cells.Item(13, 9).Value = 1.030946177668293
$ws.Cells.Item(13, 10).Value = 1.037362986162248
$ws.Cells.Item(13, 11).Value = 1.035941122810242
$ws.Cells.Item(13, 12).Value = 1.047406435300631
$ws.Cells.Item(13, 13).Value = 1.053591868704514
$ws.Cells.Item(13, 14).Value = 1.038836160589796
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.030891913073003
$ws.Cells.Item(14, 4).Value = 1.032462742024397
$ws.Cells.Item(14, 5).Value = 1.04403211917997
$ws.Cells.Item(14, 6).Value = 1.050265011288858
$ws.Cells.Item(14, 9).Value = 1.030976131400131
$ws.Cells.Item(14, 10).Value = 1.037497001152881
$ws.Cells.Item(14, 11).Value = 1.036046743959872
$ws.Cells.Item(14, 12).Value = 1.047573387547565
$ws.Cells.Item(14, 13).Value = 1.053783688332759
$ws.Cells.Item(14, 14).Value = 1.038970365897087
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.031012388928415
$ws.Cells.Item(15, 4).Value = 1.032547687100472
$ws.Cells.Item(15, 5).Value = 1.04415475086626
$ws.Cells.Item(15, 6).Value = 1.050402889376505
$ws.Cells.Item(15, 9).Value = 1.030994550556785
$ws.Cells.Item(15, 10).Value = 1.037579558362764
$ws.Cells.Item(15, 11).Value = 1.036111797502156
$ws.Cells.Item(15, 12).Value = 1.047676247196004
$ws.Cells.Item(15, 13).Value = 1.053901876795467
$ws.Cells.Item(15, 14).Value = 1.039053040347679
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.031713590063033
$ws.Cells.Item(16, 4).Value = 1.033041996081579
$ws.Cells.Item(16, 5).Value = 1.044868656403631
$ws.Cells.Item(16, 6).Value = 1.051205633738874
$ws.Cells.Item(16, 9).Value = 1.031101203332891
$ws.Cells.Item(16, 10).Value = 1.038059866044833
$ws.Cells.Item(16, 11).Value = 1.036490084186864
$ws.Cells.Item(16, 12).Value = 1.048274857016699
$ws.Cells.Item(16, 13).Value = 1.054589817092101
$ws.Cells.Item(16, 14).Value = 1.039534030121747
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.032153418190639
$ws.Cells.Item(17, 4).Value = 1.033351967875257
$ws.Cells.Item(17, 5).Value = 1.045316591265487
$ws.Cells.Item(17, 6).Value = 1.05170938167663
$ws.Cells.Item(17, 9).Value = 1.031167615254929
$ws.Cells.Item(17, 10).Value = 1.038360965104026
$ws.Cells.Item(17, 11).Value = 1.036727062685105
$ws.Cells.Item(17, 12).Value = 1.048650281403266
$ws.Cells.Item(17, 13).Value = 1.055021373378226
$ws.Cells.Item(17, 14).Value = 1.039835556776141
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.03240995689095
$ws.Cells.Item(18, 4).Value = 1.033532734884922
$ws.Cells.Item(18, 5).Value = 1.045577907968097
$ws.Cells.Item(18, 6).Value = 1.052003284633986
$ws.Cells.Item(18, 9).Value = 1.03120617613463
$ws.Cells.Item(18, 10).Value = 1.038536524539192
$ws.Cells.Item(18, 11).Value = 1.036865176324623
$ws.Cells.Item(18, 12).Value = 1.048869235987873
$ws.Cells.Item(18, 13).Value = 1.055273103891128
$ws.Cells.Item(18, 14).Value = 1.040011365525843
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.03249742905061
$ws.Cells.Item(19, 4).Value = 1.03359436594747
$ws.Cells.Item(19, 5).Value = 1.04566701768583
$ws.Cells.Item(19, 6).Value = 1.052103510767152
$ws.Cells.Item(19, 9).Value = 1.031219294532911
$ws.Cells.Item(19, 10).Value = 1.038596374501546
$ws.Cells.Item(19, 11).Value = 1.036912250532238
$ws.Cells.Item(19, 12).Value = 1.048943889782467
$ws.Cells.Item(19, 13).Value = 1.055358939399559
$ws.Cells.Item(19, 14).Value = 1.040071300482007
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.032106229384839
$ws.Cells.Item(20, 4).Value = 1.033318714382442
$ws.Cells.Item(20, 5).Value = 1.045268527563347
$ws.Cells.Item(20, 6).Value = 1.051655326511628
$ws.Cells.Item(20, 9).Value = 1.031160508095869
$ws.Cells.Item(20, 10).Value = 1.038328666912961
$ws.Cells.Item(20, 11).Value = 1.036701648711929
$ws.Cells.Item(20, 12).Value = 1.048610004419856
$ws.Cells.Item(20, 13).Value = 1.054975070340482
$ws.Cells.Item(20, 14).Value = 1.039803212717941
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.030834332374525
$ws.Cells.Item(21, 4).Value = 1.032422141401733
$ws.Cells.Item(21, 5).Value = 1.04397351088994
$ws.Cells.Item(21, 6).Value = 1.05019911785083
$ws.Cells.Item(21, 9).Value = 1.030967318340242
$ws.Cells.Item(21, 10).Value = 1.037457539964626
$ws.Cells.Item(21, 11).Value = 1.036015645983636
$ws.Cells.Item(21, 12).Value = 1.047524225344719
$ws.Cells.Item(21, 13).Value = 1.053727201806526
$ws.Cells.Item(21, 14).Value = 1.038930848669419
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.030034579975509
$ws.Cells.Item(22, 4).Value = 1.031858121319368
$ws.Cells.Item(22, 5).Value = 1.043159670307017
$ws.Cells.Item(22, 6).Value = 1.04928420963385
$ws.Cells.Item(22, 9).Value = 1.030844266173366
$ws.Cells.Item(22, 10).Value = 1.03690922373963
$ws.Cells.Item(22, 11).Value = 1.035583318957284
$ws.Cells.Item(22, 12).Value = 1.046841329243196
$ws.Cells.Item(22, 13).Value = 1.052942708228262
$ws.Cells.Item(22, 14).Value = 1.03838175377249
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.030458550040753
$ws.Cells.Item(23, 4).Value = 1.032157148339073
$ws.Cells.Item(23, 5).Value = 1.043591066202749
$ws.Cells.Item(23, 6).Value = 1.049769156983638
$ws.Cells.Item(23, 9).Value = 1.030909648793713
$ws.Cells.Item(23, 10).Value = 1.037199954032538
$ws.Cells.Item(23, 11).Value = 1.035812599502129
$ws.Cells.Item(23, 12).Value = 1.047203366907995
$ws.Cells.Item(23, 13).Value = 1.053358574913557
$ws.Cells.Item(23, 14).Value = 1.038672896935762
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.032127551991854
$ws.Cells.Item(24, 4).Value = 1.033333740310282
$ws.Cells.Item(24, 5).Value = 1.045290245345533
$ws.Cells.Item(24, 6).Value = 1.051679751480769
$ws.Cells.Item(24, 9).Value = 1.031163720059094
$ws.Cells.Item(24, 10).Value = 1.038343261280619
$ws.Cells.Item(24, 11).Value = 1.036713132538084
$ws.Cells.Item(24, 12).Value = 1.048628203929119
$ws.Cells.Item(24, 13).Value = 1.054995992655665
$ws.Cells.Item(24, 14).Value = 1.039817827811274
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.034064120921656
$ws.Cells.Item(25, 4).Value = 1.034697764468802
$ws.Cells.Item(25, 5).Value = 1.047263773089387
$ws.Cells.Item(25, 6).Value = 1.053899849542171
$ws.Cells.Item(25, 9).Value = 1.031451648597406
$ws.Cells.Item(25, 10).Value = 1.039667397270931
$ws.Cells.Item(25, 11).Value = 1.037753761799261
$ws.Cells.Item(25, 12).Value = 1.050280705519919
$ws.Cells.Item(25, 13).Value = 1.056896561876181
$ws.Cells.Item(25, 14).Value = 1.041143844226573
